# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Updates the "Periodo Mora" detail table (rows 16-25) on Hoja1 so that it
# reflects new account-statement periods (1708-1711) for the three workers,
# replacing the previous period groupings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data for the detail table (B:G), row by row.
# Columns: B=Tipo Doc (unchanged "CC"), C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row = 16; Doc = "45537634"; Nombre = "LUZ KETTY MUNARRIZ PERNETT";        Periodo = "1711"; Valor = 104880; Salario = 2622000 },
    @{ Row = 17; Doc = "45537634"; Nombre = "LUZ KETTY MUNARRIZ PERNETT";        Periodo = "1710"; Valor = 104880; Salario = 2622000 },
    @{ Row = 18; Doc = "45537634"; Nombre = "LUZ KETTY MUNARRIZ PERNETT";        Periodo = "1709"; Valor = 104880; Salario = 2622000 },
    @{ Row = 19; Doc = "45537634"; Nombre = "LUZ KETTY MUNARRIZ PERNETT";        Periodo = "1708"; Valor = 104880; Salario = 2622000 },
    @{ Row = 20; Doc = "45480236"; Nombre = "LEDA MARIA AHUMADA MOUTHON";        Periodo = "1710"; Valor = 128240; Salario = 3206000 },
    @{ Row = 21; Doc = "45480236"; Nombre = "LEDA MARIA AHUMADA MOUTHON";        Periodo = "1709"; Valor = 128240; Salario = 3206000 },
    @{ Row = 22; Doc = "45480236"; Nombre = "LEDA MARIA AHUMADA MOUTHON";        Periodo = "1708"; Valor = 128240; Salario = 3206000 },
    @{ Row = 23; Doc = "45524166"; Nombre = "DEYRA ANNGELMINA GRANADOS BELTRAN"; Periodo = "1710"; Valor = 31720;  Salario = 793000 },
    @{ Row = 24; Doc = "45524166"; Nombre = "DEYRA ANNGELMINA GRANADOS BELTRAN"; Periodo = "1709"; Valor = 31720;  Salario = 793000 },
    @{ Row = 25; Doc = "45524166"; Nombre = "DEYRA ANNGELMINA GRANADOS BELTRAN"; Periodo = "1708"; Valor = 31720;  Salario = 793000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
